$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Credentials sheet: no value changes, just selection / tab-state churn.
# ---------------------------------------------------------------------------
$wsCredentials = $wb.Worksheets.Item("Credentials")
$wsCredentials.Range("C9").Select()

# ---------------------------------------------------------------------------
# 2. Email sheet: move the selection to C4.
# ---------------------------------------------------------------------------
$wsEmail = $wb.Worksheets.Item("Email")
$wsEmail.Range("C4").Select()

# ---------------------------------------------------------------------------
# 3. ProductDetails sheet: Product/Size change from t-shirt/M to dresses/L.
# ---------------------------------------------------------------------------
$wsProductDetails = $wb.Worksheets.Item("ProductDetails")
$wsProductDetails.Range("A2").Value = "dresses"
$wsProductDetails.Range("C2").Value = "L"
$wsProductDetails.Range("C2").Select()

# ---------------------------------------------------------------------------
# 4. SearchProduct sheet: Product changes from t-shirt to dresses.
# ---------------------------------------------------------------------------
$wsSearchProduct = $wb.Worksheets.Item("SearchProduct")
$wsSearchProduct.Range("A2").Value = "dresses"
$wsSearchProduct.Range("F17").Select()

# ---------------------------------------------------------------------------
# 5. AccountCreationData sheet: update rows 2-4 with new test data, and split
#    off columns I:O (Company..MobilePhone) into a brand-new "Sheet1".
# ---------------------------------------------------------------------------
$wsAccount = $wb.Worksheets.Item("AccountCreationData")

# Capture the Company..MobilePhone block (columns I:O) before we touch it.
$company  = $wsAccount.Range("I2").Value2
$address  = $wsAccount.Range("J2").Value2
$city     = $wsAccount.Range("K2").Value2
$state    = $wsAccount.Range("L2").Value2
$zipcode  = $wsAccount.Range("M2").Value2
$country  = $wsAccount.Range("N2").Value2
$mobile   = $wsAccount.Range("O2").Value2

$hdrCompany  = $wsAccount.Range("I1").Value2
$hdrAddress  = $wsAccount.Range("J1").Value2
$hdrCity     = $wsAccount.Range("K1").Value2
$hdrState    = $wsAccount.Range("L1").Value2
$hdrZipcode  = $wsAccount.Range("M1").Value2
$hdrCountry  = $wsAccount.Range("N1").Value2
$hdrMobile   = $wsAccount.Range("O1").Value2

# New email / password / year values for rows 2-4.
$wsAccount.Range("A2").Value = "newtest31@gmail.com"
$wsAccount.Range("B2").Value = "Mr "
$wsAccount.Range("E2").Value = "hgsdtyf"
$wsAccount.Range("H2").Value = 1985

$wsAccount.Range("A3").Value = "newtest32@gmail.com"
$wsAccount.Range("B3").Value = "Mrs"
$wsAccount.Range("E3").Value = "jddjysj"
$wsAccount.Range("H3").Value = 1986

$wsAccount.Range("A4").Value = "newtest33@gmail.com"
$wsAccount.Range("B4").Value = "Mrs"
$wsAccount.Range("E4").Value = "hssuujg"
$wsAccount.Range("H4").Value = 1987

# Remove the Company..MobilePhone columns from this sheet - they move to the
# new sheet below.
$wsAccount.Range("I1:O4").Clear()

# ---------------------------------------------------------------------------
# 6. New "Sheet1" at the end of the workbook, holding the data that used to
#    live in AccountCreationData columns I:O.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastSheet)
$wsNew.Name = "Sheet1"

$wsNew.Range("A1").Value = $hdrCompany
$wsNew.Range("B1").Value = $hdrAddress
$wsNew.Range("C1").Value = $hdrCity
$wsNew.Range("D1").Value = $hdrState
$wsNew.Range("E1").Value = $hdrZipcode
$wsNew.Range("F1").Value = $hdrCountry
$wsNew.Range("G1").Value = $hdrMobile

for ($r = 2; $r -le 4; $r++) {
    $wsNew.Range("A$r").Value = $company
    $wsNew.Range("B$r").Value = $address
    $wsNew.Range("C$r").Value = $city
    $wsNew.Range("D$r").Value = $state
    $wsNew.Range("E$r").Value = $zipcode
    $wsNew.Range("F$r").Value = $country
    $wsNew.Range("G$r").Value = $mobile
}

$wsNew.Range("K11").Select()

# ---------------------------------------------------------------------------
# 7. Activate AccountCreationData last so it ends up as the selected tab,
#    matching the new activeTab index + tabSelected="1" placement.
# ---------------------------------------------------------------------------
$wsAccount.Activate()
$wsAccount.Range("K9").Select()
